# Add the new "eur_sek" column (AT) to the quarterly_averages worksheet.
# Mirrors the author's edit: a header label in AT1 (matching the style of
# the other header cells) plus one numeric value per data row (AT2:AT61).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# AT1 header - copy the formatting from the existing header cell (AS1) so
# the new column matches the rest of the header row, then set its text.
$ws.Range("AS1").Copy() | Out-Null
$ws.Range("AT1").PasteSpecial(-4122) | Out-Null
$ws.Range("AT1").Value = "eur_sek"

# AT2:AT61 - the eur_sek quarterly average values.
$ws.Range("AT2").Value = 9.9464047619047626
$ws.Range("AT3").Value = 9.631269841269841
$ws.Range("AT4").Value = 9.3803606060606057
$ws.Range("AT5").Value = 9.2138969696969699
$ws.Range("AT6").Value = 8.8641562500000006
$ws.Range("AT7").Value = 9.0152873015872999
$ws.Range("AT8").Value = 9.1450712121212128
$ws.Range("AT9").Value = 9.0910359374999992
$ws.Range("AT10").Value = 8.8528907692307701
$ws.Range("AT11").Value = 8.9133419354838708
$ws.Range("AT12").Value = 8.4354169230769234
$ws.Range("AT13").Value = 8.6230265624999998
$ws.Range("AT14").Value = 8.496509677419354
$ws.Range("AT15").Value = 8.5652365079365094
$ws.Range("AT16").Value = 8.6798000000000002
$ws.Range("AT17").Value = 8.857528125
$ws.Range("AT18").Value = 8.8569190476190478
$ws.Range("AT19").Value = 9.0516838709677412
$ws.Range("AT20").Value = 9.2051515151515151
$ws.Range("AT21").Value = 9.2717718750000007
$ws.Range("AT22").Value = 9.3799952380952387
$ws.Range("AT23").Value = 9.299512903225807
$ws.Range("AT24").Value = 9.4293303030303015
$ws.Range("AT25").Value = 9.3021646153846156
$ws.Range("AT26").Value = 9.326653225806453
$ws.Range("AT27").Value = 9.2782323076923081
$ws.Range("AT28").Value = 9.5105893939393944
$ws.Range("AT29").Value = 9.7573593750000001
$ws.Range("AT30").Value = 9.5062769230769231
$ws.Range("AT31").Value = 9.6917677419354842
$ws.Range("AT32").Value = 9.5567646153846155
$ws.Range("AT33").Value = 9.7930380952380958
$ws.Range("AT34").Value = 9.9711952380952376
$ws.Range("AT35").Value = 10.330344444444441
$ws.Range("AT36").Value = 10.40546307692308
$ws.Range("AT37").Value = 10.320371874999999
$ws.Range("AT38").Value = 10.418688888888891
$ws.Range("AT39").Value = 10.61907741935484
$ws.Range("AT40").Value = 10.662187878787879
$ws.Range("AT41").Value = 10.652359375
$ws.Range("AT42").Value = 10.668850000000001
$ws.Range("AT43").Value = 10.650670967741929
$ws.Range("AT44").Value = 10.36415606060606
$ws.Range("AT45").Value = 10.26768153846154
$ws.Range("AT46").Value = 10.12016825396825
$ws.Range("AT47").Value = 10.14136666666667
$ws.Range("AT48").Value = 10.194848484848491
$ws.Range("AT49").Value = 10.12804242424242
$ws.Range("AT50").Value = 10.4806734375
$ws.Range("AT51").Value = 10.478549206349211
$ws.Range("AT52").Value = 10.619298484848491
$ws.Range("AT53").Value = 10.937740625
$ws.Range("AT54").Value = 11.202999999999999
$ws.Range("AT55").Value = 11.469095161290319
$ws.Range("AT56").Value = 11.76409846153846
$ws.Range("AT57").Value = 11.47838253968254
$ws.Range("AT58").Value = 11.27923015873016
$ws.Range("AT59").Value = 11.50352222222222
$ws.Range("AT60").Value = 11.45124393939394
$ws.Range("AT61").Value = 11.494209375000001

# Leave the workbook with the same selection the author ended up with.
$ws.Range("AT63").Select() | Out-Null
